$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.857148
$ws.Range("H2").Value = 2.571444
$ws.Range("I2").Value = 0.04787301688248034
$ws.Range("J2").Value = 0.04787301688248034
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2315733333333333
$ws.Range("N2").Value = 0.69472
$ws.Range("O2").Value = 0.2371078251520917
$ws.Range("P2").Value = 0.2371078251520917
$ws.Range("Q2").Value = 0.19849261952
$ws.Range("R2").Value = 1.78643357568
$ws.Range("S2").Value = 0.01135106691647428
$ws.Range("T2").Value = 0.01135106691647428

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.857148
$ws.Range("H3").Value = 2.571444
$ws.Range("I3").Value = 0.04787301688248034
$ws.Range("J3").Value = 0.04787301688248034
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7200953333333334
$ws.Range("N3").Value = 2.160286
$ws.Range("O3").Value = 0.7373052671097876
$ws.Range("P3").Value = 0.7373052671097876
$ws.Range("Q3").Value = 0.6172282747760001
$ws.Range("R3").Value = 5.555054472984001
$ws.Range("S3").Value = 0.03529702749988853
$ws.Range("T3").Value = 0.03529702749988854

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.857148
$ws.Range("H4").Value = 2.571444
$ws.Range("I4").Value = 0.04787301688248034
$ws.Range("J4").Value = 0.04787301688248034
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02498966666666666
$ws.Range("N4").Value = 0.074969
$ws.Range("O4").Value = 0.02558690773812063
$ws.Range("P4").Value = 0.02558690773812063
$ws.Range("Q4").Value = 0.021419842804
$ws.Range("R4").Value = 0.192778585236
$ws.Range("S4").Value = 0.001224922466117515
$ws.Range("T4").Value = 0.001224922466117516

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.94278666666667
$ws.Range("H5").Value = 47.82836
$ws.Range("I5").Value = 0.8904288352152905
$ws.Range("J5").Value = 0.8904288352152906
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2315733333333333
$ws.Range("N5").Value = 0.69472
$ws.Range("O5").Value = 0.2371078251520917
$ws.Range("P5").Value = 0.2371078251520917
$ws.Range("Q5").Value = 3.691924251022222
$ws.Range("R5").Value = 33.2273182592
$ws.Range("S5").Value = 0.2111276445706078
$ws.Range("T5").Value = 0.2111276445706078

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl12"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.94278666666667
$ws.Range("H6").Value = 47.82836
$ws.Range("I6").Value = 0.8904288352152905
$ws.Range("J6").Value = 0.8904288352152906
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7200953333333334
$ws.Range("N6").Value = 2.160286
$ws.Range("O6").Value = 0.7373052671097876
$ws.Range("P6").Value = 0.7373052671097876
$ws.Range("Q6").Value = 11.48032627899556
$ws.Range("R6").Value = 103.32293651096
$ws.Range("S6").Value = 0.6565178701906668
$ws.Range("T6").Value = 0.6565178701906669

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl12"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.94278666666667
$ws.Range("H7").Value = 47.82836
$ws.Range("I7").Value = 0.8904288352152905
$ws.Range("J7").Value = 0.8904288352152906
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02498966666666666
$ws.Range("N7").Value = 0.074969
$ws.Range("O7").Value = 0.02558690773812063
$ws.Range("P7").Value = 0.02558690773812063
$ws.Range("Q7").Value = 0.3984049245377778
$ws.Range("R7").Value = 3.58564432084
$ws.Range("S7").Value = 0.02278332045401585
$ws.Range("T7").Value = 0.02278332045401586

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl12"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 1.104681666666667
$ws.Range("H8").Value = 3.314045
$ws.Range("I8").Value = 0.06169814790222908
$ws.Range("J8").Value = 0.06169814790222908
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2315733333333333
$ws.Range("N8").Value = 0.69472
$ws.Range("O8").Value = 0.2371078251520917
$ws.Range("P8").Value = 0.2371078251520917
$ws.Range("Q8").Value = 0.2558148158222223
$ws.Range("R8").Value = 2.3023333424
$ws.Range("S8").Value = 0.01462911366500963
$ws.Range("T8").Value = 0.01462911366500963

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl12"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 1.104681666666667
$ws.Range("H9").Value = 3.314045
$ws.Range("I9").Value = 0.06169814790222908
$ws.Range("J9").Value = 0.06169814790222908
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7200953333333334
$ws.Range("N9").Value = 2.160286
$ws.Range("O9").Value = 0.7373052671097876
$ws.Range("P9").Value = 0.7373052671097876
$ws.Range("Q9").Value = 0.7954761129855558
$ws.Range("R9").Value = 7.159285016870001
$ws.Range("S9").Value = 0.0454903694192322
$ws.Range("T9").Value = 0.0454903694192322

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ccl12"
$ws.Range("C10").Value = "Ackr4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 1.104681666666667
$ws.Range("H10").Value = 3.314045
$ws.Range("I10").Value = 0.06169814790222908
$ws.Range("J10").Value = 0.06169814790222908
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02498966666666666
$ws.Range("N10").Value = 0.074969
$ws.Range("O10").Value = 0.02558690773812063
$ws.Range("P10").Value = 0.02558690773812063
$ws.Range("Q10").Value = 0.02760562662277778
$ws.Range("R10").Value = 0.248450639605
$ws.Range("S10").Value = 0.001578664817987256
$ws.Range("T10").Value = 0.001578664817987257
